$wb = $excel.ActiveWorkbook

# --- ALC sheet (index 1) ---
$ws = $wb.Worksheets.Item(1)
# Row 17
$ws.Range("H17").Value = 3274.6667
$ws.Range("J17").Value = 3274.6667
$ws.Range("L17").Value = 9824.000100000001
$ws.Range("N17").Value = -10160.0001

# Row 29
$ws.Range("H29").Value = 3780
$ws.Range("I29").Value = 3749
$ws.Range("J29").Value = 3787.75
$ws.Range("K29").Value = 11247
$ws.Range("L29").Value = 11363.25
$ws.Range("M29").Value = -10966
$ws.Range("N29").Value = -11925.25

# Row 38
$ws.Range("H38").Value = 1484.2
$ws.Range("I38").Value = 230.25
$ws.Range("J38").Value = 6500
$ws.Range("K38").Value = 690.75
$ws.Range("L38").Value = 19500
$ws.Range("M38").Value = -318.75
$ws.Range("N38").Value = -20244

# Row 92
$ws.Range("H92").Value = 291.75
$ws.Range("I92").Value = 267.5
$ws.Range("J92").Value = 364.5
$ws.Range("K92").Value = 267.5
$ws.Range("L92").Value = 364.5
$ws.Range("M92").Value = 980.5
$ws.Range("N92").Value = -2860.5

# Row 94
$ws.Range("H94").Value = 3328.6
$ws.Range("I94").Value = 3328.6
$ws.Range("K94").Value = 3328.6
$ws.Range("M94").Value = -2877.6

# Row 96
$ws.Range("H96").Value = 901.875
$ws.Range("I96").Value = 745.1429000000001
$ws.Range("K96").Value = 2235.4287
$ws.Range("M96").Value = -862.4287000000004

# Row 99
$ws.Range("H99").Value = 326.25
$ws.Range("J99").Value = 297
$ws.Range("L99").Value = 891
$ws.Range("N99").Value = -3887

# Row 100
$ws.Range("H100").Value = 1752.0834
$ws.Range("I100").Value = 2192.4
$ws.Range("J100").Value = 1437.5714
$ws.Range("K100").Value = 2192.4
$ws.Range("L100").Value = 1437.5714
$ws.Range("M100").Value = -1651.4
$ws.Range("N100").Value = -2519.5714

# Row 112
$ws.Range("H112").Value = 1170.8649
$ws.Range("J112").Value = 1521.75
$ws.Range("L112").Value = 4565.25
$ws.Range("N112").Value = -6781.25

# Row 116
$ws.Range("H116").Value = 5678.1113
$ws.Range("I116").Value = 5418.5
$ws.Range("J116").Value = 6197.3335
$ws.Range("K116").Value = 5418.5
$ws.Range("L116").Value = 6197.3335
$ws.Range("M116").Value = -1976.5
$ws.Range("N116").Value = -13081.3335

# Row 118
$ws.Range("H118").Value = 1591.6666
$ws.Range("I118").Value = 987.5
$ws.Range("K118").Value = 2962.5
$ws.Range("M118").Value = -1305.5

# Row 127
$ws.Range("H127").Value = 4398.5
$ws.Range("I127").Value = 797
$ws.Range("J127").Value = 8000
$ws.Range("K127").Value = 2391
$ws.Range("L127").Value = 24000
$ws.Range("M127").Value = 2569
$ws.Range("N127").Value = -33920

# Row 138
$ws.Range("H138").Value = 2739.2034
$ws.Range("I138").Value = 1983.5454
$ws.Range("J138").Value = 2912.375
$ws.Range("K138").Value = 5950.6362
$ws.Range("L138").Value = 8737.125
$ws.Range("M138").Value = -810.6361999999999
$ws.Range("N138").Value = -19017.125

# --- ARM sheet (index 2) ---
$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Range("H32").Value = 6365.3
$ws.Range("I32").Value = 6365.3
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 6365.3
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -6078.3
$ws.Range("N32").ClearContents()

# Row 37
$ws.Range("H37").Value = 31798.8
$ws.Range("I37").Value = 6999
$ws.Range("K37").Value = 6999
$ws.Range("M37").Value = -6726

# Row 74
$ws.Range("H74").Value = 2499
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 2499
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 2499
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -4247

# Row 77
$ws.Range("H77").Value = 2499
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 2499
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 12495
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -21231

# Row 122
$ws.Range("H122").Value = 1878.2858
$ws.Range("I122").Value = 1891.5
$ws.Range("K122").Value = 5674.5
$ws.Range("M122").Value = -3224.5

# Row 132
$ws.Range("H132").Value = 3236.9167
$ws.Range("I132").Value = 2640.8333
$ws.Range("K132").Value = 7922.499899999999
$ws.Range("M132").Value = -5392.499899999999

# --- BSM sheet (index 3) ---
$ws = $wb.Worksheets.Item(3)
# Row 35
$ws.Range("H35").Value = 41268.5
$ws.Range("J35").Value = 41268.5
$ws.Range("L35").Value = 41268.5
$ws.Range("N35").Value = -41888.5

# Row 80
$ws.Range("H80").Value = 708
$ws.Range("J80").Value = 1060
$ws.Range("L80").Value = 1060
$ws.Range("N80").Value = -3056

# Row 83
$ws.Range("H83").Value = 708
$ws.Range("J83").Value = 1060
$ws.Range("L83").Value = 5300
$ws.Range("N83").Value = -15284

# Row 99
$ws.Range("H99").Value = 1804.8889
$ws.Range("I99").Value = 1813.4286
$ws.Range("J99").Value = 1775
$ws.Range("K99").Value = 1813.4286
$ws.Range("L99").Value = 1775
$ws.Range("M99").Value = -315.4286
$ws.Range("N99").Value = -4771

# Row 105
$ws.Range("H105").Value = 3739.2
$ws.Range("I105").Value = 2799
$ws.Range("J105").Value = 7500
$ws.Range("K105").Value = 2799
$ws.Range("L105").Value = 7500
$ws.Range("M105").Value = -1052
$ws.Range("N105").Value = -10994

# --- CRP sheet (index 4) ---
$ws = $wb.Worksheets.Item(4)
# Row 17
$ws.Range("H17").Value = 15999.5
$ws.Range("J17").Value = 17000
$ws.Range("L17").Value = 17000
$ws.Range("N17").Value = -17348

# Row 25
$ws.Range("H25").Value = 20013
$ws.Range("J25").Value = 20013
$ws.Range("L25").Value = 20013
$ws.Range("N25").Value = -20361

# Row 31
$ws.Range("H31").Value = 2535.875
$ws.Range("J31").Value = 2014
$ws.Range("L31").Value = 2014
$ws.Range("N31").Value = -2604

# Row 34
$ws.Range("H34").Value = 2535.875
$ws.Range("J34").Value = 2014
$ws.Range("L34").Value = 2014
$ws.Range("N34").Value = -2418

# Row 41
$ws.Range("H41").Value = 39965
$ws.Range("J41").Value = 39965
$ws.Range("L41").Value = 39965
$ws.Range("N41").Value = -40821

# Row 107
$ws.Range("H107").Value = 974.875
$ws.Range("I107").Value = 525.75
$ws.Range("K107").Value = 525.75
$ws.Range("M107").Value = 1394.25

# Row 122
$ws.Range("H122").Value = 4155.625
$ws.Range("I122").Value = 3849
$ws.Range("K122").Value = 11547
$ws.Range("M122").Value = -9097

# Row 134
$ws.Range("H134").Value = 3674.25
$ws.Range("I134").Value = 3999
$ws.Range("K134").Value = 11997
$ws.Range("M134").Value = -9462

# --- CUL sheet (index 5) ---
$ws = $wb.Worksheets.Item(5)
# Row 34
$ws.Range("H34").Value = 3799.8
$ws.Range("J34").Value = 5433.3335
$ws.Range("L34").Value = 16300.0005
$ws.Range("N34").Value = -16468.0005

# Row 39
$ws.Range("H39").Value = 10000
$ws.Range("J39").Value = 10000
$ws.Range("L39").Value = 30000
$ws.Range("N39").Value = -30588

# Row 55
$ws.Range("H55").Value = 5000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 5000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 15000
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -15354

# Row 56
$ws.Range("H56").Value = 17851.8
$ws.Range("I56").Value = 17851.8
$ws.Range("K56").Value = 17851.8
$ws.Range("M56").Value = -17321.8

# --- GSM sheet (index 6) ---
$ws = $wb.Worksheets.Item(6)
# Row 43
$ws.Range("H43").Value = 31495
$ws.Range("J43").Value = 31495
$ws.Range("L43").Value = 31495
$ws.Range("N43").Value = -31797

# Row 46
$ws.Range("H46").Value = 19961.375
$ws.Range("J46").Value = 19961.375
$ws.Range("L46").Value = 19961.375
$ws.Range("N46").Value = -20273.375

# Row 122
$ws.Range("H122").Value = 1683.1428
$ws.Range("I122").Value = 1714
$ws.Range("J122").Value = 1498
$ws.Range("K122").Value = 5142
$ws.Range("L122").Value = 4494
$ws.Range("M122").Value = -2692
$ws.Range("N122").Value = -9394

# --- LTW sheet (index 7) ---
$ws = $wb.Worksheets.Item(7)
# Row 40
$ws.Range("H40").Value = 4040.6667
$ws.Range("I40").Value = 4040.6667
$ws.Range("K40").Value = 4040.6667
$ws.Range("M40").Value = -3904.6667

# Row 132
$ws.Range("H132").Value = 2902.8333
$ws.Range("I132").Value = 1981.7778
$ws.Range("K132").Value = 5945.3334
$ws.Range("M132").Value = -3415.3334

# --- WVR sheet (index 8) ---
$ws = $wb.Worksheets.Item(8)
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
